# Apply updated TPM-derived values to the Slitrk3-Ptprs sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.02439366666666666
$ws.Range("H2").Value = 0.073181
$ws.Range("M2").Value = 3.556762333333333
$ws.Range("N2").Value = 10.670287
$ws.Range("O2").Value = 0.04280930450251701
$ws.Range("P2").Value = 0.04280930450251701
$ws.Range("Q2").Value = 0.08676247477188888
$ws.Range("R2").Value = 0.7808622729469999
$ws.Range("S2").Value = 0.04280930450251701
$ws.Range("T2").Value = 0.04280930450251701

# Row 3
$ws.Range("G3").Value = 0.02439366666666666
$ws.Range("H3").Value = 0.073181
$ws.Range("O3").Value = 0.5686906263805706
$ws.Range("P3").Value = 0.5686906263805704
$ws.Range("Q3").Value = 1.152576681582222
$ws.Range("R3").Value = 10.37319013424
$ws.Range("S3").Value = 0.5686906263805706
$ws.Range("T3").Value = 0.5686906263805704

# Row 4
$ws.Range("G4").Value = 0.02439366666666666
$ws.Range("H4").Value = 0.073181
$ws.Range("M4").Value = 24.53173066666666
$ws.Range("N4").Value = 73.595192
$ws.Range("O4").Value = 0.2952646900921413
$ws.Range("P4").Value = 0.2952646900921412
$ws.Range("Q4").Value = 0.598418860639111
$ws.Range("R4").Value = 5.385769745752
$ws.Range("S4").Value = 0.2952646900921413
$ws.Range("T4").Value = 0.2952646900921412

# Row 5
$ws.Range("G5").Value = 0.02439366666666666
$ws.Range("H5").Value = 0.073181
$ws.Range("M5").Value = 7.746355333333334
$ws.Range("N5").Value = 23.239066
$ws.Range("O5").Value = 0.09323537902477132
$ws.Range("P5").Value = 0.0932353790247713
$ws.Range("Q5").Value = 0.1889620098828889
$ws.Range("R5").Value = 1.700658088946
$ws.Range("S5").Value = 0.09323537902477132
$ws.Range("T5").Value = 0.0932353790247713
